$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "93.346.11"
Set-TextValue "E2" "  -5.34%  "
Set-TextValue "D3" "3.362.30"
Set-TextValue "E3" "  -3.33%  "
Set-TextValue "D5" "233.77"
Set-TextValue "E5" "  -8.56%  "
Set-TextValue "D6" "625.42"
Set-TextValue "E6" "  -7.17%  "
Set-TextValue "D7" "1.38"
Set-TextValue "E7" "  -7.97%  "
Set-TextValue "D8" "0.388"
Set-TextValue "E8" "  -10.88%  "
Set-TextValue "E9" "  +0.18%  "
Set-TextValue "D10" "0.939"
Set-TextValue "E10" "  -11.53%  "
Set-TextValue "D11" "3.364.58"
Set-TextValue "E11" "  -3.18%  "
Set-TextValue "E12" "  -7.90%  "
Set-TextValue "D13" "40.26"
Set-TextValue "E13" "  -13.61%  "
Set-TextValue "E14" "  -3.52%  "
Set-TextValue "D15" "93.173.77"
Set-TextValue "E15" "  -5.42%  "
Set-TextValue "D16" "3.990.32"
Set-TextValue "E16" "  -3.35%  "
Set-TextValue "D17" "0.0000243"
Set-TextValue "E17" "  -6.99%  "
Set-TextValue "D18" "7.99"
Set-TextValue "E18" "  -12.40%  "
Set-TextValue "D19" "3.362.02"
Set-TextValue "E19" "  -3.14%  "
Set-TextValue "D20" "16.86"
Set-TextValue "E20" "  -9.29%  "
Set-TextValue "D21" "10.88"
Set-TextValue "E21" "  -7.51%  "
Set-TextValue "D22" "490.80"
Set-TextValue "E22" "  -6.27%  "
Set-TextValue "D23" "0.449"
Set-TextValue "E23" "  -16.80%  "
Set-TextValue "D24" "3.13"
Set-TextValue "E24" "  -9.83%  "
Set-TextValue "D25" "0.0000185"
Set-TextValue "E25" "  -9.63%  "
Set-TextValue "D26" "6.17"
Set-TextValue "E26" "  -10.38%  "
Set-TextValue "D27" "89.87"
Set-TextValue "E27" "  -8.78%  "
Set-TextValue "D28" "3.548.18"
Set-TextValue "E28" "  -2.96%  "
Set-TextValue "D29" "11.46"
Set-TextValue "E29" "  -10.11%  "
Set-TextValue "D30" "11.30"
Set-TextValue "E30" "  -8.10%  "
Set-TextValue "E31" "  +0.06%  "
Set-TextValue "D32" "2.64"
Set-TextValue "E32" "  -10.02%  "
Set-TextValue "E33" "  -10.94%  "
Set-TextValue "E34" "  -0.38%  "
Set-TextValue "D35" "0.171"
Set-TextValue "E35" "  -11.21%  "
Set-TextValue "D36" "28.63"
Set-TextValue "E36" "  -5.24%  "
Set-TextValue "D37" "0.527"
Set-TextValue "E37" "  -9.61%  "
Set-TextValue "D38" "7.46"
Set-TextValue "E38" "  -8.74%  "
Set-TextValue "D39" "521.73"
Set-TextValue "E39" "  -2.86%  "
Set-TextValue "E40" "  -0.04%  "
Set-TextValue "E41" "  -9.60%  "
Set-TextValue "D42" "0.147"
Set-TextValue "E42" "  -5.83%  "
Set-TextValue "D43" "0.876"
Set-TextValue "E43" "  -2.02%  "
Set-TextValue "D44" "24.03"
Set-TextValue "E44" "  -1.68%  "
Set-TextValue "D45" "3.59"
Set-TextValue "E45" "  -2.31%  "
Set-TextValue "E46" "  -7.65%  "
Set-TextValue "D47" "5.48"
Set-TextValue "E47" "  -5.40%  "
Set-TextValue "D48" "2.14"
Set-TextValue "E48" "  -4.91%  "
Set-TextValue "D49" "0.0391"
Set-TextValue "E49" "  -11.17%  "
Set-TextValue "D50" "52.32"
Set-TextValue "E50" "  -6.63%  "
Set-TextValue "D51" "7.95"
Set-TextValue "E51" "  -9.65%  "
